$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) for rows 2-12: 45212 -> 45221 (serial date value)
for ($row = 2; $row -le 12; $row++) {
    $ws.Cells.Item($row, 3).Value = 45221
}
